# Word COM-interop script implementing the diff:
#   1. 'arızalı'  -> 'faulty'     (inside the ClassNames cell array literal)
#   2. 'arızasız' -> 'faultless'  (inside the ClassNames cell array literal)
#   3. Delete the two explanatory paragraphs
#         faulty=‘arızalı’
#         faultless= ‘arızasız’
#      (together with one of the blank paragraphs that framed them, so the
#      run of blank paragraphs collapses from two separate ones around the
#      removed text into a single contiguous block of blank paragraphs).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) 'arızalı' -> 'faulty'
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("'arızalı'", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    # Assigning .Text directly (instead of using Find's ReplaceWith) keeps the
    # literal straight apostrophes instead of letting AutoCorrect turn them
    # into curly quotes.
    $rng.Text = "'faulty"
    $rng.Collapse(0)
    $rng.InsertAfter("'")
}

# ---------------------------------------------------------------------
# 2) 'arızasız' -> 'faultless'
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("'arızasız'", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "'faultless"
    $rng.Collapse(0)
    $rng.InsertAfter("'")
}

# ---------------------------------------------------------------------
# 3) Remove the two "faulty=..." / "faultless=..." explanatory paragraphs
#    entirely (together with the paragraph mark of the blank paragraph
#    that precedes them), leaving the rest of the blank paragraphs intact.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("faulty=", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $paraCount = $d.Paragraphs.Count
    $startParaIndex = 0
    for ($i = 1; $i -le $paraCount; $i++) {
        $p = $d.Paragraphs($i)
        if (($rng.Start -ge $p.Range.Start) -and ($rng.Start -lt $p.Range.End)) {
            $startParaIndex = $i
        }
    }

    if ($startParaIndex -gt 0) {
        $firstPara = $d.Paragraphs($startParaIndex)
        $secondPara = $d.Paragraphs($startParaIndex + 1)
        $delRange = $d.Range($firstPara.Range.Start, $secondPara.Range.End)
        $delRange.Delete()
    }
}
